$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "gray" row's complementary color changes from #6A3524 to #81402C.
$ws.Range("F3").Value2 = "#81402C"

# New "vector" column: a header plus a single combined color-scale string
# (styled in Helvetica) on the "gray" row.
$ws.Range("G1").Value2 = "vector"
$ws.Range("G3").Value2 = "#15343D, #284046, #3E4D4E, #5F5B4E, #856949, #AB7743, #D2863C, #FB9637, #E57630, #D0562B, #BC3626, #A61922, #8D0422, #70002A, #520036"
$ws.Range("G3").Font.Name = "Helvetica"

# Leave the selection where the author last left it.
[void]$ws.Range("G17").Select()
